$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 147 (shifts old rows 147-211 down to 148-212)
$ws.Rows(147).Insert()

# Populate the newly inserted row 147 with the new weekly record.
$ws.Range("A147").Value = 4
$ws.Range("B147").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C147").Value = "Los Lagos"
$ws.Range("D147").Value = 44596
$ws.Range("E147").Value = 10
$ws.Range("F147").Value = "Fruta"
$ws.Range("G147").Value = 100104
$ws.Range("H147").Value = "Frutos de pepita"
$ws.Range("I147").Value = 100104005
$ws.Range("J147").Value = "Pera"
$ws.Range("K147").Value = "Packham's Triumph"
$ws.Range("L147").Value = "Primera"
$ws.Range("M147").Value = 400
$ws.Range("N147").Value = 14000
$ws.Range("O147").Value = 15000
$ws.Range("P147").Value = 14500
$ws.Range("Q147").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R147").Value = "Región de O'Higgins"
$ws.Range("S147").Value = 967
$ws.Range("T147").Value = 15
